$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Processes sheet: swap the "description" (B) and "type" (E) columns.
#    Before: key | description | fuel | product | type
#    After:  key | type        | fuel | product | description
# ---------------------------------------------------------------------------
$wsProcesses = $wb.Worksheets.Item("Processes")
for ($r = 1; $r -le 7; $r++) {
    $bVal = $wsProcesses.Cells.Item($r, 2).Value()
    $eVal = $wsProcesses.Cells.Item($r, 5).Value()
    $wsProcesses.Cells.Item($r, 2).Value = $eVal
    $wsProcesses.Cells.Item($r, 5).Value = $bVal
}

# Resize columns B and E to roughly match the new content widths, and shrink
# the unused, vestigial column G formatting down to F.
$wsProcesses.Columns.Item(2).ColumnWidth = 12.451822916666666
$wsProcesses.Columns.Item(5).ColumnWidth = 13.307291666666666
$wsProcesses.Columns.Item(6).ColumnWidth = 9.166666666666666

$wsProcesses.Range("B7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Update the "cgam_processes" defined name range (E7 -> D7).
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Processes!cgam_processes") {
        $n.RefersTo = "=Processes!`$A`$1:`$D`$7"
    }
}

# ---------------------------------------------------------------------------
# 3. WasteDefinition sheet: add the new "recycle" column (C).
# ---------------------------------------------------------------------------
$wsWasteDef = $wb.Worksheets.Item("WasteDefinition")
$wsWasteDef.Range("C1").Value = "recycle"
$wsWasteDef.Range("C2").Value = 0
$wsWasteDef.Range("B1").Copy() | Out-Null
$wsWasteDef.Range("C1").PasteSpecial(-4122) | Out-Null
$wsWasteDef.Range("C1").Select() | Out-Null
$wsWasteDef.Columns.Item(3).Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. WasteAllocation sheet: move the selection to C1.
# ---------------------------------------------------------------------------
$wsWasteAlloc = $wb.Worksheets.Item("WasteAllocation")
$wsWasteAlloc.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Make WasteDefinition the active/selected sheet (was Format before).
# ---------------------------------------------------------------------------
$wsWasteDef.Activate() | Out-Null
$wsWasteDef.Columns.Item(3).Select() | Out-Null
